# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Swap "Valor Mora" values between the first period row (2108, row 16)
# and the last period row (2102, row 22).
$ws.Range("F16").Value = 36341
$ws.Range("F22").Value = 35129
